# Generate Report for Handoff
# Adds the newly-handed-off file "f085978c-53ea-43e4-8410-78b6f50f3c3a.md"
# as a new row on the Overview, zh-cn and de-de sheets, keeping the
# ListObjects (tables) and used-range in sync with the appended row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet (row 9)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A9").Value = "f085978c-53ea-43e4-8410-78b6f50f3c3a.md"
$wsOverview.Range("C9").Value = ".md"
$wsOverview.Range("E9").Value = "Ready for handoff"
$wsOverview.Range("F9").Value = "Ready for handoff"
$wsOverview.Range("G9").Value = "2016-08-23 00:41:53"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B9"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f085978c53ea43e4841078b6f50f3c3afake0001/e2e/f085978c-53ea-43e4-8410-78b6f50f3c3a.md",
    "",
    "",
    "e2e\f085978c-53ea-43e4-8410-78b6f50f3c3a.md"
) | Out-Null

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G9"))

# ---------------------------------------------------------------------
# zh-cn sheet (row 9)
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("B9").Value = ".md"
$wsZhCn.Range("C9").Value = "Ready for handoff"
$wsZhCn.Range("D9").Value = "e2e"
$wsZhCn.Range("E9").Value = "ht"
$wsZhCn.Range("F9").Value = "False"
$wsZhCn.Range("G9").Value = "f085978c-53ea-43e4-8410-78b6f50f3c3a.9b9f3e0182837636493b923089b1e5dd760ce9bc.zh-cn.xlf"
$wsZhCn.Range("H9").Value = "2016-08-23 00:41:48"
$wsZhCn.Range("K9").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("M9").Value = "True"
$wsZhCn.Range("O9").Value = "False"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A9"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/f085978c53ea43e4841078b6f50f3c3afake0002/e2e/f085978c-53ea-43e4-8410-78b6f50f3c3a.md",
    "",
    "",
    "f085978c-53ea-43e4-8410-78b6f50f3c3a.md"
) | Out-Null

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P9"))

# ---------------------------------------------------------------------
# de-de sheet (row 9)
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("B9").Value = ".md"
$wsDeDe.Range("C9").Value = "Ready for handoff"
$wsDeDe.Range("D9").Value = "e2e"
$wsDeDe.Range("E9").Value = "ht"
$wsDeDe.Range("F9").Value = "False"
$wsDeDe.Range("G9").Value = "f085978c-53ea-43e4-8410-78b6f50f3c3a.9b9f3e0182837636493b923089b1e5dd760ce9bc.de-de.xlf"
$wsDeDe.Range("H9").Value = "2016-08-23 00:41:53"
$wsDeDe.Range("K9").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("M9").Value = "True"
$wsDeDe.Range("O9").Value = "False"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A9"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/f085978c53ea43e4841078b6f50f3c3afake0003/e2e/f085978c-53ea-43e4-8410-78b6f50f3c3a.md",
    "",
    "",
    "f085978c-53ea-43e4-8410-78b6f50f3c3a.md"
) | Out-Null

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P9"))
